# Add a new "Lua" column (G) to the worksheet, mirroring the existing
# Matlab/Python/C/Cpp/Js/Java columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(1, 7).Value = "Lua"

# New Lua data for rows 2-11
$luaValues = @(3.377, 8.458, 8.947, 8.65, 8.844, 9.336, 9.465, 9.613, 3.008, 8.78)
for ($i = 0; $i -lt $luaValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $luaValues[$i]
}

# G2 carries the wrapped-text look used elsewhere in the sheet (e.g. E5)
$ws.Cells.Item(2, 7).WrapText = $true

# Move/record the active selection on the newly edited cell, as in the source file
$ws.Range("G10").Select()
